# Generate Report for Handoff
# Replaces the localization-status report data: the previous handoff/handback
# cycle (a3673701.../ec16c88b...) is superseded by a fresh "Ready for handoff"
# cycle (7945b2a5.../ffff04920d58...), clearing the now-stale handback columns.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ws1.Range("A2").Value = "7945b2a5-f081-4612-9729-8eb5f78851c0.md"
$ws1.Range("B2").Value = "e2e\7945b2a5-f081-4612-9729-8eb5f78851c0.md"
$ws1.Range("E2").Value = "Ready for handoff"
$ws1.Range("F2").Value = "Ready for handoff"
$ws1.Range("G2").Value = "2016-08-27 09:01:53"

$ws1.Range("A3").Value = "ffff04920d58-557b-4ec2-85f7-a6c77be4f265.md"
$ws1.Range("B3").Value = "e2e\ffff04920d58-557b-4ec2-85f7-a6c77be4f265.md"
$ws1.Range("E3").Value = "Ready for handoff"
$ws1.Range("F3").Value = "Ready for handoff"
$ws1.Range("G3").Value = "2016-08-27 09:01:53"

foreach ($h in $ws1.Hyperlinks) {
    if ($h.Range.Address() -eq '$B$2') {
        $h.TextToDisplay = "e2e\7945b2a5-f081-4612-9729-8eb5f78851c0.md"
    } elseif ($h.Range.Address() -eq '$B$3') {
        $h.TextToDisplay = "e2e\ffff04920d58-557b-4ec2-85f7-a6c77be4f265.md"
    }
}

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws2.Range("A2").Value = "7945b2a5-f081-4612-9729-8eb5f78851c0.md"
$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("G2").Value = "7945b2a5-f081-4612-9729-8eb5f78851c0.62e8cad1dcd6bf6aa0348605d9897ce8587f75a2.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-08-27 09:01:49"
$ws2.Range("I2").Value = ""
$ws2.Range("J2").Value = ""
$ws2.Range("K2").Value = "0001-01-01 00:00:00"

$ws2.Range("A3").Value = "ffff04920d58-557b-4ec2-85f7-a6c77be4f265.md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("F3").Value = "True"
$ws2.Range("G3").Value = "7945b2a5-f081-4612-9729-8eb5f78851c0.62e8cad1dcd6bf6aa0348605d9897ce8587f75a2.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-27 09:01:49"
$ws2.Range("I3").Value = ""
$ws2.Range("J3").Value = ""
$ws2.Range("K3").Value = "0001-01-01 00:00:00"

$toDelete2 = @()
foreach ($h in $ws2.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$2') {
        $h.TextToDisplay = "7945b2a5-f081-4612-9729-8eb5f78851c0.md"
    } elseif ($h.Range.Address() -eq '$A$3') {
        $h.TextToDisplay = "ffff04920d58-557b-4ec2-85f7-a6c77be4f265.md"
    } elseif ($h.Range.Address() -eq '$I$2') {
        $toDelete2 += $h
    } elseif ($h.Range.Address() -eq '$I$3') {
        $toDelete2 += $h
    }
}
foreach ($h in $toDelete2) {
    $h.Delete()
}

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws3.Range("A2").Value = "7945b2a5-f081-4612-9729-8eb5f78851c0.md"
$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("G2").Value = "7945b2a5-f081-4612-9729-8eb5f78851c0.62e8cad1dcd6bf6aa0348605d9897ce8587f75a2.de-de.xlf"
$ws3.Range("H2").Value = "2016-08-27 09:01:53"
$ws3.Range("I2").Value = ""
$ws3.Range("J2").Value = ""
$ws3.Range("K2").Value = "0001-01-01 00:00:00"

$ws3.Range("A3").Value = "ffff04920d58-557b-4ec2-85f7-a6c77be4f265.md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("F3").Value = "True"
$ws3.Range("G3").Value = "7945b2a5-f081-4612-9729-8eb5f78851c0.62e8cad1dcd6bf6aa0348605d9897ce8587f75a2.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-27 09:01:53"
$ws3.Range("I3").Value = ""
$ws3.Range("J3").Value = ""
$ws3.Range("K3").Value = "0001-01-01 00:00:00"

$toDelete3 = @()
foreach ($h in $ws3.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$2') {
        $h.TextToDisplay = "7945b2a5-f081-4612-9729-8eb5f78851c0.md"
    } elseif ($h.Range.Address() -eq '$A$3') {
        $h.TextToDisplay = "ffff04920d58-557b-4ec2-85f7-a6c77be4f265.md"
    } elseif ($h.Range.Address() -eq '$I$2') {
        $toDelete3 += $h
    } elseif ($h.Range.Address() -eq '$I$3') {
        $toDelete3 += $h
    }
}
foreach ($h in $toDelete3) {
    $h.Delete()
}

# ---------------------------------------------------------------------------
# Column width autofit - the shortened cell content (e.g. "Ready for handoff"
# replacing "Handed back: in sync with en-US", and now-empty I/J columns)
# narrows these columns.
# ---------------------------------------------------------------------------
$ws1.Columns("E:F").EntireColumn.AutoFit()
$ws2.Columns("C").EntireColumn.AutoFit()
$ws2.Columns("I:J").EntireColumn.AutoFit()
$ws3.Columns("C").EntireColumn.AutoFit()
$ws3.Columns("I:J").EntireColumn.AutoFit()
